# Apply cell value updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / Link text columns (plain text, no numeric coercion risk)
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

# Price / Volume(1h) columns -- force text format so numeric-looking
# strings (e.g. "2.570", "0.005958") keep their exact textual digits
# instead of being reinterpreted as numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '315.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.61%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '42.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-5.30%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.188'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.42%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08042'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.04%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.374'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.53%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.728'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-11.24%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9275'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-4.83%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1123'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.39%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1844'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-3.04%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09213'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.10%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04566'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.10%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.354'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-15.23%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.46%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001274'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.01%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04187'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.95%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005958'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3.77%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.354'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.34%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.570'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.49%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3385'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.67%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1385'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.68%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2635'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.13%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001247'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004242'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.14%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001228'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.72%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002992'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.30%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02528'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-8.46%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05390'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.54%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.008015'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.95%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1388'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.70%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007614'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.06%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002071'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.95%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008417'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.28%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3138'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-10.32%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006762'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.59%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000753'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.30%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003399'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-2.70%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004118'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '16.50%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002109'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.30%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002009'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.30%'
